$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144, shifting existing rows 144-182 down to 145-183.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new weekly price record.
$ws.Cells.Item(144, 1).Value  = 11
$ws.Cells.Item(144, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(144, 3).Value  = "Bíobío"
$ws.Cells.Item(144, 4).Value  = 44736
$ws.Cells.Item(144, 5).Value  = 8
$ws.Cells.Item(144, 6).Value  = 100112003
$ws.Cells.Item(144, 7).Value  = "Ajo"
$ws.Cells.Item(144, 8).Value  = "Chino"
$ws.Cells.Item(144, 9).Value  = "Primera"
$ws.Cells.Item(144, 10).Value = 400
$ws.Cells.Item(144, 11).Value = 17000
$ws.Cells.Item(144, 12).Value = 18000
$ws.Cells.Item(144, 13).Value = 17500
$ws.Cells.Item(144, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(144, 15).Value = "China"
$ws.Cells.Item(144, 16).Value = 1750
$ws.Cells.Item(144, 17).Value = 10
$ws.Cells.Item(144, 18).Value = "Hortaliza"
